$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 876874
$ws.Range("C20").Value = 3149718
$ws.Range("C26").Value = 4559224
$ws.Range("C30").Value = 214083
$ws.Range("C31").Value = 42102199
$ws.Range("C38").Value = 153682
$ws.Range("C40").Value = 270113
$ws.Range("C43").Value = 201428
$ws.Range("C46").Value = 628842
$ws.Range("C49").Value = 7979486
$ws.Range("C51").Value = 153318
$ws.Range("C53").Value = 99586
$ws.Range("C56").Value = 36188
$ws.Range("C58").Value = 1408781
$ws.Range("C60").Value = 663886
$ws.Range("C62").Value = 1181486
$ws.Range("C64").Value = 181441
$ws.Range("C66").Value = 17067899
$ws.Range("C68").Value = 451209
$ws.Range("C71").Value = 814720
$ws.Range("C72").Value = 375494
$ws.Range("C73").Value = 284268
$ws.Range("C75").Value = 328272
$ws.Range("C79").Value = 931206
$ws.Range("C82").Value = 1641801
$ws.Range("C84").Value = 120523
$ws.Range("C86").Value = 232091
$ws.Range("C88").Value = 119732
$ws.Range("C90").Value = 386157
$ws.Range("C92").Value = 3392220
$ws.Range("C94").Value = 64446712
$ws.Range("C99").Value = 467539
$ws.Range("C101").Value = 643908
$ws.Range("C104").Value = 2275929
$ws.Range("C106").Value = 660209
$ws.Range("C109").Value = 551428
$ws.Range("C112").Value = 1643384
$ws.Range("C116").Value = 655645
$ws.Range("C118").Value = 1555897
$ws.Range("C120").Value = 513387
$ws.Range("C122").Value = 817865
$ws.Range("C124").Value = 1592530
$ws.Range("C128").Value = 1017847
$ws.Range("C129").Value = 1702043
$ws.Range("C138").Value = 176102
$ws.Range("C140").Value = 226540
$ws.Range("C144").Value = 137937
$ws.Range("C145").Value = 263073
$ws.Range("C146").Value = 165265
$ws.Range("C157").Value = 163657
$ws.Range("C165").Value = 4516809
$ws.Range("C188").Value = 1405486
$ws.Range("C206").Value = 2355097
$ws.Range("C232").Value = 939746
$ws.Range("C245").Value = 578898
$ws.Range("C252").Value = 969603
$ws.Range("C265").Value = 450149
$ws.Range("C274").Value = 16316774
$ws.Range("C290").Value = 299341
$ws.Range("C291").Value = 212126
$ws.Range("C297").Value = 133758
